$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.035.97"
$ws.Range("D3").Value = "1.830.21"
$ws.Range("E3").Value = "  -0.14%  "
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "0.9990"
$r.Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "241.24"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -0.34%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "0.6229"
$r.Style = "Normal"
$ws.Range("E6").Value = "  -5.37%  "
$ws.Range("E7").Value = "  +0.00%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.07553"
$r.Style = "Normal"
$ws.Range("E8").Value = "  +2.16%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "44.64"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +6.72%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.2909"
$r.Style = "Normal"
$ws.Range("E10").Value = "  -0.59%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "22.77"
$r.Style = "Normal"
$ws.Range("E11").Value = "  -0.37%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.07635"
$r.Style = "Normal"
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").Value = "1.831.38"
$ws.Range("E13").Value = "  -0.03%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "4.957"
$r.Style = "Normal"
$ws.Range("E14").Value = "  -0.75%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "0.6651"
$r.Style = "Normal"
$ws.Range("E15").Value = "  -0.05%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "82.30"
$r.Style = "Normal"
$ws.Range("E16").Value = "  -0.62%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "0.000009101"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +8.35%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "5.991"
$r.Style = "Normal"
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("D19").Value = "29.035.90"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").Value = "2.084.70"
$ws.Range("E20").Value = "  +0.26%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "224.70"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -1.05%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "12.34"
$r.Style = "Normal"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +0.75%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "1.001"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +0.58%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "8.380"
$r.Style = "Normal"
$ws.Range("E27").Value = "  -2.57%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "0.1356"
$r.Style = "Normal"
$ws.Range("E28").Value = "  -2.27%  "
$ws.Range("E29").Value = "  -0.48%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "1.493"
$r.Style = "Normal"
$ws.Range("E30").Value = "  -1.65%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "4.052"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("E32").Value = "  +0.92%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "4.025"
$r.Style = "Normal"
$ws.Range("E33").Value = "  -0.42%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "0.05205"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -1.22%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "1.836"
$r.Style = "Normal"
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("E36").Value = "  +1.34%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.7322"
$r.Style = "Normal"
$ws.Range("E37").Value = "  -1.10%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "2.604"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("D39").Value = "1.285.76"
$ws.Range("E39").Value = "  -1.24%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "2.761"
$r.Style = "Normal"
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("E41").Value = "  -0.59%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "6.374"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +7.59%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.8900"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -3.66%  "
$ws.Range("E44").Value = "  +0.04%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "101.38"
$r.Style = "Normal"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").Value = "1.981.43"
$ws.Range("E46").Value = "  +0.51%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.5115"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -0.53%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "8.805"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +0.28%  "
